# Alterada ordem dos slides
# Swap the order of slide 7 ("Orientação a Objetos x Relacional" table)
# and slide 8 ("JAVA PERSISTENCE API") so the JPA slide now comes before
# the table slide.

$p = $ppt.ActivePresentation

# Move slide 7 to position 8 — this pushes the current slide 8 up to
# position 7, effectively swapping the two slides.
$s = $p.Slides.Item(7)
$s.MoveTo(8)
